$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each content-bearing row of the worksheet table, 5 columns of "a÷b=" problems.
# Mapping derived in row/column order to avoid ambiguity from duplicate values.
$rows = @(1, 5, 9, 13, 17)

$values = @(
    @("87÷5=", "23÷7=", "43÷4=", "83÷7=", "51÷9="),
    @("36÷4=", "27÷9=", "34÷9=", "48÷2=", "89÷6="),
    @("91÷6=", "52÷6=", "60÷2=", "74÷7=", "47÷9="),
    @("73÷5=", "27÷6=", "60÷2=", "18÷3=", "84÷5="),
    @("12÷5=", "58÷5=", "83÷8=", "64÷3=", "14÷3=")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    for ($c = 1; $c -le 5; $c++) {
        $t.Cell($r, $c).Range.Text = $values[$i][$c - 1]
    }
}

Write-Host "Updated table cells"
